$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H5").Value = "2\04\2022"
$ws.Range("H5").Select()
